# Add a new data row (row 19) to the "Artfynd" sheet, mirroring the
# structure of the existing rows (e.g. row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

# --- Numeric columns -------------------------------------------------
$ws.Cells.Item($row, 1).Value  = 112209223          # A  Id
$ws.Cells.Item($row, 2).Value  = 89802               # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value  = 5420                # E  TaxonId
$ws.Cells.Item($row, 17).Value = 753585.5470142091   # Q  Ost
$ws.Cells.Item($row, 18).Value = 7090378.824306269   # R  Nord
$ws.Cells.Item($row, 19).Value = 1                   # S  Noggrannhet

# --- Plain text columns (unambiguous strings - no auto-conversion risk)
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"                  # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value  = "LC"                          # D  Rödlistade
$ws.Cells.Item($row, 6).Value  = "Grovticka"                   # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Phaeolus schweinitzii"       # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Fr.) Pat."                  # H  Auktor
$ws.Cells.Item($row, 16).Value = "Umedalen, Umeå, Vb"          # P  Lokalnamn
$ws.Cells.Item($row, 20).Value = "Västerbotten"                # T  Län
$ws.Cells.Item($row, 21).Value = "Umeå"                        # U  Kommun
$ws.Cells.Item($row, 22).Value = "Västerbotten"                # V  Provins
$ws.Cells.Item($row, 23).Value = "Umeå socken"                 # W  Församling
$ws.Cells.Item($row, 49).Value = "anitha lundberg"              # AW Rapportör
$ws.Cells.Item($row, 50).Value = "anitha lundberg"              # AX Observatörer

# --- Text columns that LOOK numeric/date - force text via a leading
#     apostrophe so Excel doesn't reinterpret the literal as a number
#     or a date/time serial. ----------------------------------------
$ws.Cells.Item($row, 9).Value  = "'5"       # I  Antal
$ws.Cells.Item($row, 25).Value = "'2023-09-18"   # Y  Startdatum
$ws.Cells.Item($row, 26).Value = "'00:00"        # Z  Starttid
$ws.Cells.Item($row, 27).Value = "'2023-09-18"   # AA Slutdatum
$ws.Cells.Item($row, 28).Value = "'00:00"        # AB Sluttid

# --- Boolean columns ---------------------------------------------------
$ws.Cells.Item($row, 30).Value = $false   # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false   # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false   # AG Ospontan

# --- Empty (but present) cells, matching the source row's blank cells -
# A lone apostrophe is Excel's "empty, forced-text" literal: it yields a
# present Text cell whose value is "" (same as the sibling rows' blank
# inlineStr cells), instead of clearing/removing the cell entirely.
$ws.Cells.Item($row, 10).Value = "'"   # J  Enhet
$ws.Cells.Item($row, 11).Value = "'"   # K  Ålder-Stadium
$ws.Cells.Item($row, 14).Value = "'"   # N  Metod
$ws.Cells.Item($row, 32).Value = "'"   # AF Bestämningsmetod
$ws.Cells.Item($row, 46).Value = "'"   # AT Bestämningsår
$ws.Cells.Item($row, 51).Value = "'"   # AY Projektnamn
